$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update confidence-interval text cells (stored as inline strings)
$ws.Range("L2").Value = "[60.18295375703367, 67.75345931151467]"
$ws.Range("T2").Value = "[47.43204292614425, 52.580776373443925]"
$ws.Range("L3").Value = "[59.61403140573342, 67.77575799733069]"
$ws.Range("P3").Value = "[1.1761317842268104, 1.3019212798660407]"
$ws.Range("T3").Value = "[46.7822916248994, 51.44341833735159]"

# Update numeric cells
$ws.Range("X3").Value = 20.176576576577
$ws.Range("Y3").Value = 20.68608608608652
